# The underlying species-observation records effectively rotated one
# position through rows 75-78 (75<-77, 76<-78, 77<-76, 78<-75), while the
# row numbers / cell addresses stayed fixed. Apply that as direct
# per-cell value updates instead of physically moving rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 75 (becomes old row 77's "Lunglav" record) ---
$ws.Cells.Item(75, 1).Value = 111573439          # A75 Id
$ws.Cells.Item(75, 2).Value = 78578              # B75 Taxonsorteringsordning
$ws.Cells.Item(75, 4).Value = "NT"               # D75 Rödlistade
$ws.Cells.Item(75, 5).Value = 6458               # E75 TaxonId
$ws.Cells.Item(75, 6).Value = "Lunglav"          # F75 Artnamn
$ws.Cells.Item(75, 7).Value = "Lobaria pulmonaria"  # G75 Vetenskapligt namn
$ws.Cells.Item(75, 8).Value = "(L.) Hoffm."       # H75 Auktor
$ws.Cells.Item(75, 29).Value = "Även på en gran i närheten."  # AC75 Publik kommentar

# --- Row 76 (becomes old row 78's "Knärot" record) ---
$ws.Cells.Item(76, 1).Value = 111573395          # A76 Id
$ws.Cells.Item(76, 2).Value = 96348              # B76 Taxonsorteringsordning
$ws.Cells.Item(76, 4).Value = "VU"               # D76 Rödlistade
$ws.Cells.Item(76, 5).Value = 220787             # E76 TaxonId
$ws.Cells.Item(76, 6).Value = "Knärot"           # F76 Artnamn
$ws.Cells.Item(76, 7).Value = "Goodyera repens"  # G76 Vetenskapligt namn
$ws.Cells.Item(76, 8).Value = "(L.) R. Br."       # H76 Auktor
$ws.Cells.Item(76, 9).Value = ""                 # I76 Antal (cleared)
$ws.Cells.Item(76, 10).Value = ""                # J76 Enhet (cleared)
$ws.Cells.Item(76, 11).Value = "överblommad"     # K76 Ålder-Stadium

# --- Row 77 (becomes old row 76's "Skogsfru" record) ---
$ws.Cells.Item(77, 1).Value = 111573403          # A77 Id
$ws.Cells.Item(77, 2).Value = 96346              # B77 Taxonsorteringsordning
$ws.Cells.Item(77, 5).Value = 620                # E77 TaxonId
$ws.Cells.Item(77, 6).Value = "Skogsfru"         # F77 Artnamn
$ws.Cells.Item(77, 7).Value = "Epipogium aphyllum"  # G77 Vetenskapligt namn
$ws.Cells.Item(77, 8).Value = "Sw."               # H77 Auktor
$ws.Cells.Item(77, 9).Value = "1"                # I77 Antal
$ws.Cells.Item(77, 10).Value = "stjälkar/strån/skott"  # J77 Enhet
$ws.Cells.Item(77, 11).Value = "blomning"        # K77 Ålder-Stadium
$ws.Cells.Item(77, 17).Value = 491981.2275731571 # Q77 Ost
$ws.Cells.Item(77, 18).Value = 7015311.94031445  # R77 Nord
$ws.Cells.Item(77, 29).Value = ""                # AC77 Publik kommentar (cleared)

# --- Row 78 (becomes old row 75's "Luddlav" record) ---
$ws.Cells.Item(78, 1).Value = 111573805          # A78 Id
$ws.Cells.Item(78, 2).Value = 78612              # B78 Taxonsorteringsordning
$ws.Cells.Item(78, 4).Value = "LC"               # D78 Rödlistade
$ws.Cells.Item(78, 5).Value = 6464               # E78 TaxonId
$ws.Cells.Item(78, 6).Value = "Luddlav"          # F78 Artnamn
$ws.Cells.Item(78, 7).Value = "Nephroma resupinatum"  # G78 Vetenskapligt namn
$ws.Cells.Item(78, 8).Value = "(L.) Ach."         # H78 Auktor
$ws.Cells.Item(78, 11).Value = ""                # K78 Ålder-Stadium (cleared)
$ws.Cells.Item(78, 17).Value = 492215.3225248906 # Q78 Ost
$ws.Cells.Item(78, 18).Value = 7015165.030750753 # R78 Nord
$ws.Cells.Item(78, 29).Value = "På sälg med lunglav."  # AC78 Publik kommentar
